# Fruta / hortaliza, semanal
# Insert a new week's worth of records (rows 527-528: Primera/Segunda for date 45124)
# right after current row 526, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 527 (existing rows 527.. shift down to 529..)
$ws.Range("A527:A528").EntireRow.Insert()

# --- New row 527: Primera, 2023-07-17 (serial 45124) ---
$ws.Cells.Item(527, 1).Value = 4
$ws.Cells.Item(527, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(527, 3).Value = "Los Lagos"
$ws.Cells.Item(527, 4).Value = 45124
$ws.Cells.Item(527, 5).Value = 10
$ws.Cells.Item(527, 6).Value = "Fruta"
$ws.Cells.Item(527, 7).Value = 100101
$ws.Cells.Item(527, 8).Value = "Berries"
$ws.Cells.Item(527, 9).Value = 100101007
$ws.Cells.Item(527, 10).Value = "Kiwi"
$ws.Cells.Item(527, 11).Value = "Hayward"
$ws.Cells.Item(527, 12).Value = "Primera"
$ws.Cells.Item(527, 13).Value = 200
$ws.Cells.Item(527, 14).Value = 16000
$ws.Cells.Item(527, 15).Value = 16000
$ws.Cells.Item(527, 16).Value = 16000
$ws.Cells.Item(527, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(527, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(527, 19).Value = 1067
$ws.Cells.Item(527, 20).Value = 15

# --- New row 528: Segunda, 2023-07-17 (serial 45124) ---
$ws.Cells.Item(528, 1).Value = 4
$ws.Cells.Item(528, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(528, 3).Value = "Los Lagos"
$ws.Cells.Item(528, 4).Value = 45124
$ws.Cells.Item(528, 5).Value = 10
$ws.Cells.Item(528, 6).Value = "Fruta"
$ws.Cells.Item(528, 7).Value = 100101
$ws.Cells.Item(528, 8).Value = "Berries"
$ws.Cells.Item(528, 9).Value = 100101007
$ws.Cells.Item(528, 10).Value = "Kiwi"
$ws.Cells.Item(528, 11).Value = "Hayward"
$ws.Cells.Item(528, 12).Value = "Segunda"
$ws.Cells.Item(528, 13).Value = 200
$ws.Cells.Item(528, 14).Value = 13000
$ws.Cells.Item(528, 15).Value = 13000
$ws.Cells.Item(528, 16).Value = 13000
$ws.Cells.Item(528, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(528, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(528, 19).Value = 867
$ws.Cells.Item(528, 20).Value = 15
